$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 1394
$ws.Range("A13").Value = 1395
$ws.Range("A14").Value = 1396
$ws.Range("A15").Value = 1599
$ws.Range("A16").Value = 1600
$ws.Range("A17").Value = 1605
$ws.Range("A18").Value = 1606
$ws.Range("A19").Value = 1909
$ws.Range("A20").Value = 1910
$ws.Range("A21").Value = 1911
$ws.Range("A22").Value = 2526
$ws.Range("A23").Value = 2527
$ws.Range("A24").Value = 2863
$ws.Range("A25").Value = 2864
$ws.Range("A26").Value = 2865
$ws.Range("A27").Value = 4093
$ws.Range("A28").Value = 4094
$ws.Range("A29").Value = 4095
$ws.Range("A30").Value = 4189
$ws.Range("A31").Value = 4190
$ws.Range("A32").Value = 4675
$ws.Range("A33").Value = 4676
$ws.Range("A34").Value = 4677
$ws.Range("A35").Value = 4773
$ws.Range("A36").Value = 4774
$ws.Range("A37").Value = 5231
$ws.Range("A38").Value = 5232
$ws.Range("A39").Value = 5415
$ws.Range("A40").Value = 5416
$ws.Range("A41").Value = 5480
$ws.Range("A42").Value = 5481
$ws.Range("A43").Value = 5671
$ws.Range("A44").Value = 5672
$ws.Range("A45").Value = 5864
$ws.Range("A46").Value = 5865
$ws.Range("A47").Value = 5996
$ws.Range("A48").Value = 5997
$ws.Range("A49").Value = 6033
$ws.Range("A50").Value = 6034
$ws.Range("A51").Value = 6035
$ws.Range("A52").Value = 6236
$ws.Range("A53").Value = 6237
$ws.Range("A54").Value = 6508
$ws.Range("A55").Value = 6509
$ws.Range("A56").Value = 6606
$ws.Range("A57").Value = 6607
$ws.Range("A58").Value = 7328
$ws.Range("A59").Value = 7329
$ws.Range("A60").Value = 7820
$ws.Range("A61").Value = 7821
